$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = "aa"
$ws.Range("J8").Value = "Agree/Accept"
$ws.Range("I25").Value = "sd"
$ws.Range("J25").Value = "Statement-non-opinion"
$ws.Range("I38").Value = "sd"
$ws.Range("J38").Value = "Statement-non-opinion"
$ws.Range("I39").Value = "sd"
$ws.Range("J39").Value = "Statement-non-opinion"
$ws.Range("I46").Value = "b"
$ws.Range("J46").Value = "Acknowledge (Backchannel)"
$ws.Range("I56").Value = "sd"
$ws.Range("J56").Value = "Statement-non-opinion"
$ws.Range("I67").Value = "sd"
$ws.Range("J67").Value = "Statement-non-opinion"
$ws.Range("I77").Value = "sv"
$ws.Range("J77").Value = "Statement-opinion"
$ws.Range("I84").Value = "sd"
$ws.Range("J84").Value = "Statement-non-opinion"
$ws.Range("I90").Value = "sv"
$ws.Range("J90").Value = "Statement-opinion"
$ws.Range("I94").Value = "aa"
$ws.Range("J94").Value = "Agree/Accept"
$ws.Range("I100").Value = "sd"
$ws.Range("J100").Value = "Statement-non-opinion"
$ws.Range("I105").Value = "ba"
$ws.Range("J105").Value = "Appreciation"
$ws.Range("I111").Value = "sd"
$ws.Range("J111").Value = "Statement-non-opinion"
$ws.Range("I112").Value = "sv"
$ws.Range("J112").Value = "Statement-opinion"
$ws.Range("I113").Value = "ba"
$ws.Range("J113").Value = "Appreciation"
$ws.Range("I127").Value = "sd"
$ws.Range("J127").Value = "Statement-non-opinion"
$ws.Range("I141").Value = "sd"
$ws.Range("J141").Value = "Statement-non-opinion"
$ws.Range("I142").Value = "sv"
$ws.Range("J142").Value = "Statement-opinion"
$ws.Range("I163").Value = "%"
$ws.Range("J163").Value = "Uninterpretable"
$ws.Range("I168").Value = "sv"
$ws.Range("J168").Value = "Statement-opinion"
$ws.Range("I188").Value = "b"
$ws.Range("J188").Value = "Acknowledge (Backchannel)"
$ws.Range("I205").Value = "aa"
$ws.Range("J205").Value = "Agree/Accept"
$ws.Range("I208").Value = "sd"
$ws.Range("J208").Value = "Statement-non-opinion"
$ws.Range("I214").Value = "b"
$ws.Range("J214").Value = "Acknowledge (Backchannel)"
$ws.Range("I215").Value = "ba"
$ws.Range("J215").Value = "Appreciation"
$ws.Range("I217").Value = "sv"
$ws.Range("J217").Value = "Statement-opinion"
$ws.Range("I229").Value = "aa"
$ws.Range("J229").Value = "Agree/Accept"
$ws.Range("I230").Value = "sv"
$ws.Range("J230").Value = "Statement-opinion"
$ws.Range("I232").Value = "sv"
$ws.Range("J232").Value = "Statement-opinion"
$ws.Range("I242").Value = "ba"
$ws.Range("J242").Value = "Appreciation"
$ws.Range("I275").Value = "sd"
$ws.Range("J275").Value = "Statement-non-opinion"
$ws.Range("I297").Value = "b"
$ws.Range("J297").Value = "Acknowledge (Backchannel)"
$ws.Range("I299").Value = "b"
$ws.Range("J299").Value = "Acknowledge (Backchannel)"
$ws.Range("I302").Value = "%"
$ws.Range("J302").Value = "Uninterpretable"
$ws.Range("I313").Value = "aa"
$ws.Range("J313").Value = "Agree/Accept"
$ws.Range("I342").Value = "sv"
$ws.Range("J342").Value = "Statement-opinion"
$ws.Range("I349").Value = "qy"
$ws.Range("J349").Value = "Yes-No-Question"
$ws.Range("I358").Value = "sd"
$ws.Range("J358").Value = "Statement-non-opinion"
$ws.Range("I359").Value = "b"
$ws.Range("J359").Value = "Acknowledge (Backchannel)"
$ws.Range("I362").Value = "sd"
$ws.Range("J362").Value = "Statement-non-opinion"
$ws.Range("I367").Value = "sd"
$ws.Range("J367").Value = "Statement-non-opinion"
$ws.Range("I368").Value = "sd"
$ws.Range("J368").Value = "Statement-non-opinion"
$ws.Range("I401").Value = "aa"
$ws.Range("J401").Value = "Agree/Accept"
$ws.Range("I409").Value = "ba"
$ws.Range("J409").Value = "Appreciation"
